# Auto-generated edit script: applies the cell-level value changes described
# in the commit diff for Sheets/Bahamut_Profits.xlsx (workbook tabs ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Each sheet tab holds one FFXIV-leve profit
# table; only specific rows H..N market-price/profit columns changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 265.95456
$ws.Range("I28").Value = 262.94116
$ws.Range("J28").Value = 276.2
$ws.Range("K28").Value = 262.94116
$ws.Range("L28").Value = 276.2
$ws.Range("M28").Value = 222.05884
$ws.Range("N28").Value = -1246.2
$ws.Range("H33").Value = 126.38461
$ws.Range("I33").Value = 120.25
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 120.25
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = 108.75
$ws.Range("N33").Value = -658
$ws.Range("H53").Value = 347.5
$ws.Range("I53").Value = 92
$ws.Range("J53").Value = 603
$ws.Range("K53").Value = 92
$ws.Range("L53").Value = 603
$ws.Range("M53").Value = 545
$ws.Range("N53").Value = -1877
$ws.Range("H62").Value = 106890.4
$ws.Range("I62").Value = 171500.67
$ws.Range("J62").Value = 9975
$ws.Range("K62").Value = 171500.67
$ws.Range("L62").Value = 9975
$ws.Range("M62").Value = -170876.67
$ws.Range("N62").Value = -11223
$ws.Range("H65").Value = 106890.4
$ws.Range("I65").Value = 171500.67
$ws.Range("J65").Value = 9975
$ws.Range("K65").Value = 857503.3500000001
$ws.Range("L65").Value = 49875
$ws.Range("M65").Value = -854383.3500000001
$ws.Range("N65").Value = -56115
$ws.Range("H116").Value = 4336.875
$ws.Range("I116").Value = 4115.8335
$ws.Range("K116").Value = 4115.8335
$ws.Range("M116").Value = -673.8334999999997
$ws.Range("H123").Value = 20206.654
$ws.Range("J123").Value = 20206.654
$ws.Range("L123").Value = 20206.654
$ws.Range("N123").Value = -30006.654
$ws.Range("H128").Value = 42653.332
$ws.Range("J128").Value = 42653.332
$ws.Range("L128").Value = 42653.332
$ws.Range("N128").Value = -52613.332
$ws.Range("H132").Value = 2649.919
$ws.Range("I132").Value = 2644.9333
$ws.Range("J132").Value = 2671.2856
$ws.Range("K132").Value = 7934.7999
$ws.Range("L132").Value = 8013.8568
$ws.Range("M132").Value = -5404.7999
$ws.Range("N132").Value = -13073.8568
$ws.Range("H137").Value = 1120.091
$ws.Range("I137").Value = 790.8125
$ws.Range("J137").Value = 1998.1666
$ws.Range("K137").Value = 2372.4375
$ws.Range("L137").Value = 5994.4998
$ws.Range("M137").Value = 177.5625
$ws.Range("N137").Value = -11094.4998
$ws.Range("H141").Value = 2994.8076
$ws.Range("I141").Value = 1591.5625
$ws.Range("K141").Value = 4774.6875
$ws.Range("M141").Value = 405.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 2321
$ws.Range("J46").Value = 2321
$ws.Range("L46").Value = 2321
$ws.Range("N46").Value = -2959
$ws.Range("H74").Value = 838.1212
$ws.Range("I74").Value = 787.2857
$ws.Range("J74").Value = 1122.8
$ws.Range("K74").Value = 787.2857
$ws.Range("L74").Value = 1122.8
$ws.Range("M74").Value = 86.71429999999998
$ws.Range("N74").Value = -2870.8
$ws.Range("H77").Value = 838.1212
$ws.Range("I77").Value = 787.2857
$ws.Range("J77").Value = 1122.8
$ws.Range("K77").Value = 3936.4285
$ws.Range("L77").Value = 5614
$ws.Range("M77").Value = 431.5715
$ws.Range("N77").Value = -14350
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 3201.5
$ws.Range("I132").Value = 1404
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4212
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1682
$ws.Range("N132").Value = -20057
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 14871.1
$ws.Range("I107").Value = 1844.4286
$ws.Range("J107").Value = 45266.668
$ws.Range("K107").Value = 1844.4286
$ws.Range("L107").Value = 45266.668
$ws.Range("M107").Value = 75.57140000000004
$ws.Range("N107").Value = -49106.668
$ws.Range("H128").Value = 2002939
$ws.Range("I128").Value = 2002939
$ws.Range("K128").Value = 6008817
$ws.Range("M128").Value = -6006327
$ws.Range("H132").Value = 89111.125
$ws.Range("J132").Value = 89111.125
$ws.Range("L132").Value = 89111.125
$ws.Range("N132").Value = -99231.125
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H140").Value = 72523.55
$ws.Range("J140").Value = 72523.55
$ws.Range("L140").Value = 72523.55
$ws.Range("N140").Value = -82883.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2416.9119
$ws.Range("I31").Value = 2453.1667
$ws.Range("K31").Value = 2453.1667
$ws.Range("M31").Value = -2158.1667
$ws.Range("H34").Value = 2416.9119
$ws.Range("I34").Value = 2453.1667
$ws.Range("K34").Value = 2453.1667
$ws.Range("M34").Value = -2251.1667
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H62").Value = 6675
$ws.Range("I62").Value = 7633.3335
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 7633.3335
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -7009.3335
$ws.Range("N62").Value = -5048
$ws.Range("H65").Value = 6675
$ws.Range("I65").Value = 7633.3335
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 38166.6675
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -35046.6675
$ws.Range("N65").Value = -25240
$ws.Range("H116").Value = 32750
$ws.Range("J116").Value = 32750
$ws.Range("L116").Value = 32750
$ws.Range("N116").Value = -41928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 900
$ws.Range("J68").Value = 1100
$ws.Range("L68").Value = 3300
$ws.Range("N68").Value = -4922
$ws.Range("H70").Value = 4307.9
$ws.Range("I70").Value = 3324.75
$ws.Range("J70").Value = 4963.3335
$ws.Range("K70").Value = 9974.25
$ws.Range("L70").Value = 14890.0005
$ws.Range("M70").Value = -9659.25
$ws.Range("N70").Value = -15520.0005
$ws.Range("H71").Value = 900
$ws.Range("J71").Value = 1100
$ws.Range("L71").Value = 9900
$ws.Range("N71").Value = -18012
$ws.Range("H73").Value = 4307.9
$ws.Range("I73").Value = 3324.75
$ws.Range("J73").Value = 4963.3335
$ws.Range("K73").Value = 9974.25
$ws.Range("L73").Value = 14890.0005
$ws.Range("M73").Value = -8882.25
$ws.Range("N73").Value = -17074.0005
$ws.Range("H80").Value = 7285.7144
$ws.Range("I80").Value = 3666.6667
$ws.Range("K80").Value = 11000.0001
$ws.Range("M80").Value = -10064.0001
$ws.Range("H83").Value = 7285.7144
$ws.Range("I83").Value = 3666.6667
$ws.Range("K83").Value = 33000.0003
$ws.Range("M83").Value = -28320.0003
$ws.Range("H117").Value = 2742.1428
$ws.Range("J117").Value = 3096.75
$ws.Range("L117").Value = 9290.25
$ws.Range("N117").Value = -16174.25
$ws.Range("H129").Value = 1837.7894
$ws.Range("I129").Value = 1736.6666
$ws.Range("J129").Value = 1884.4615
$ws.Range("K129").Value = 5209.9998
$ws.Range("L129").Value = 5653.3845
$ws.Range("M129").Value = -209.9997999999996
$ws.Range("N129").Value = -15653.3845
$ws.Range("H141").Value = 6740
$ws.Range("I141").Value = 2726
$ws.Range("K141").Value = 8178
$ws.Range("M141").Value = -2998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5386.1665
$ws.Range("I70").Value = 4102.6665
$ws.Range("J70").Value = 6669.6665
$ws.Range("K70").Value = 4102.6665
$ws.Range("L70").Value = 6669.6665
$ws.Range("M70").Value = -3832.6665
$ws.Range("N70").Value = -7209.6665
$ws.Range("H73").Value = 5386.1665
$ws.Range("I73").Value = 4102.6665
$ws.Range("J73").Value = 6669.6665
$ws.Range("K73").Value = 4102.6665
$ws.Range("L73").Value = 6669.6665
$ws.Range("M73").Value = -3166.6665
$ws.Range("N73").Value = -8541.666499999999
$ws.Range("H113").Value = 16333.333
$ws.Range("I113").Value = 23500
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 23500
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -21330
$ws.Range("N113").Value = -6340
$ws.Range("H126").Value = 2502.5454
$ws.Range("I126").Value = 4971.3335
$ws.Range("J126").Value = 1576.75
$ws.Range("K126").Value = 14914.0005
$ws.Range("L126").Value = 4730.25
$ws.Range("M126").Value = -12444.0005
$ws.Range("N126").Value = -9670.25
$ws.Range("H132").Value = 3300.8965
$ws.Range("I132").Value = 3320.9333
$ws.Range("J132").Value = 3279.4285
$ws.Range("K132").Value = 9962.7999
$ws.Range("L132").Value = 9838.2855
$ws.Range("M132").Value = -7432.7999
$ws.Range("N132").Value = -14898.2855
$ws.Range("H133").Value = 32978
$ws.Range("J133").Value = 32978
$ws.Range("L133").Value = 32978
$ws.Range("N133").Value = -43098

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 499.77777
$ws.Range("J22").Value = 462.5
$ws.Range("L22").Value = 462.5
$ws.Range("N22").Value = -1052.5
$ws.Range("H27").Value = 499.77777
$ws.Range("J27").Value = 462.5
$ws.Range("L27").Value = 462.5
$ws.Range("N27").Value = -676.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 76553.75
$ws.Range("J133").Value = 76553.75
$ws.Range("L133").Value = 76553.75
$ws.Range("N133").Value = -86673.75
